$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("order2")
$ws2 = $wb.Worksheets.Item("order3")
$ws3 = $wb.Worksheets.Item("functionalOutput")

# ---------------------------------------------------------------------------
# order2 (sheet1): P-value column was text like "0.0174 *  " -- replace with
# the bare numeric p-value and format it as 0.0000
# ---------------------------------------------------------------------------
$ws1.Range("E2").Value = 0.0174
$ws1.Range("E3").Value = 0.00000000513
$ws1.Range("E4").Value = 0.0221
$ws1.Range("E5").Value = 0.0000163
$ws1.Range("E6").Value = 0.0000461
$ws1.Range("E7").Value = 0.00000171
$ws1.Range("E2:E7").NumberFormat = "0.0000"

# ---------------------------------------------------------------------------
# order3 (sheet2): same treatment, plus it already had one numeric p-value
# (row 6) and one "< 2e-16 ***" entry (row 7) that becomes a plain 0.
# ---------------------------------------------------------------------------
$ws2.Range("E2").Value = 0.00415
$ws2.Range("E3").Value = 0.00000000227
$ws2.Range("E4").Value = 0.00000154
$ws2.Range("E5").Value = 0.0000000123
$ws2.Range("E6").Value = 0.29055
$ws2.Range("E7").Value = 0
$ws2.Range("E8").Value = 0.000000000348
$ws2.Range("E9").Value = 0.0000336
$ws2.Range("E10").Value = 0.0000123
$ws2.Range("E11").Value = 0.000000000000239
$ws2.Range("E12").Value = 0.0000000757
$ws2.Range("E13").Value = 0.00000637
$ws2.Range("E14").Value = 0.000000833
$ws2.Range("E2:E14").NumberFormat = "0.0000"

# column sizing on order3 -- column A (coefficient names) and column E
# (p-values) were widened/best-fit after the edit
$ws2.Columns.Item(1).AutoFit() | Out-Null
$ws2.Columns.Item(5).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$ws2.Range("E8").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("A1:E15").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("J10").Select() | Out-Null
